$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-08-30"

# Update the column header label (shared string used by I1)
$ws.Range("I1").Value = "2022 (through 08-30)"

# Update August (row 9) and Total (row 14) values in the 2022 column (I)
$ws.Range("I9").Value = 163
$ws.Range("I14").Value = 1134
